$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Artfynd")

# The edit swaps the species-identifying data between row 13 and row 14,
# while leaving the other columns (location, dates, observers, etc.) untouched.
$cols = @("A", "B", "D", "E", "F", "G", "H", "Q", "R")

foreach ($col in $cols) {
    $cell13 = $ws.Range($col + "13")
    $cell14 = $ws.Range($col + "14")
    $tmp = $cell13.Value()
    $cell13.Value = $cell14.Value()
    $cell14.Value = $tmp
}
